$wb = $excel.ActiveWorkbook

$bValues = @{}
$bValues[1] = @(63,74,83,92,101,109,117,124,130,130,130,130,130,130,130,130,130,130,130,130,130)
$bValues[2] = @(62,72,82,91,99,108,115,122,129,130,130,130,130,130,130,130,130,130,130,130,130)
$bValues[3] = @(60,70,80,89,98,106,113,120,127,130,130,130,130,130,130,130,130,130,130,130,130)
$bValues[4] = @(58,69,78,88,96,104,112,119,125,130,130,130,130,130,130,130,130,130,130,130,130)
$bValues[5] = @(57,67,77,86,95,103,110,117,124,130,130,130,130,130,130,130,130,130,130,130,130)
$bValues[6] = @(55,66,75,85,93,101,109,116,122,128,130,130,130,130,130,130,130,130,130,130,130)
$bValues[7] = @(53,64,74,83,92,100,107,114,121,127,130,130,130,130,130,130,130,130,130,130,130)
$bValues[8] = @(51,62,72,81,90,98,105,112,119,125,130,130,130,130,130,130,130,130,130,130,130)
$bValues[9] = @(49,60,70,79,88,96,103,110,117,123,128,130,130,130,130,130,130,130,130,130,130)
$bValues[10] = @(48,58,68,78,86,94,102,109,115,121,126,130,130,130,130,130,130,130,130,130,130)
$bValues[11] = @(46,56,67,76,85,93,100,107,113,119,125,130,130,130,130,130,130,130,130,130,130)
$bValues[12] = @(44,55,65,74,83,91,99,106,112,118,123,128,130,130,130,130,130,130,130,130,130)
$bValues[13] = @(43,53,64,73,82,90,97,104,110,116,122,127,130,130,130,130,130,130,130,130,130)
$bValues[14] = @(41,52,62,72,81,89,96,103,109,115,120,125,130,130,130,130,130,130,130,130,130)
$bValues[15] = @(40,51,61,71,79,87,95,102,108,114,119,124,129,130,130,130,130,130,130,130,130)
$bValues[16] = @(40,50,60,69,78,86,94,101,107,113,118,123,128,130,130,130,130,130,130,130,130)
$bValues[17] = @(40,48,59,68,77,85,92,99,105,111,116,121,126,130,130,130,130,130,130,130,130)
$bValues[18] = @(40,47,57,67,75,83,91,97,104,109,114,119,124,128,130,130,130,130,130,130,130)
$bValues[19] = @(40,45,56,65,74,82,89,96,102,107,112,117,121,126,129,130,130,130,130,130,130)
$bValues[20] = @(40,45,55,64,73,81,88,94,100,105,110,115,119,123,127,130,130,130,130,130,130)

for ($i = 1; $i -le 20; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Rows.Item(2).Insert()
    $ws.Range("A2:B2").ClearFormats()
    $vals = $bValues[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $row = 2 + $j
        $ws.Cells.Item($row, 1).Value = $j
        $ws.Cells.Item($row, 2).Value = $vals[$j]
    }
}
